$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.142321392454221
$ws.Range("D2").Value = 0.1263816796427193
$ws.Range("E2").Value = 0.1492635381992038
$ws.Range("F2").Value = 2.039106289188112
$ws.Range("G2").Value = 1.373840508530236
$ws.Range("H2").Value = 1.263369999138547
$ws.Range("J2").Value = 0.2019439887166499
$ws.Range("K2").Value = 1.059503800523885
$ws.Range("B3").Value = 0.1330042420125181
$ws.Range("D3").Value = 0.1230043056577728
$ws.Range("E3").Value = 0.1458335090149347
$ws.Range("F3").Value = 2.032384371547067
$ws.Range("G3").Value = 1.367861571895276
$ws.Range("H3").Value = 1.267010523859824
$ws.Range("J3").Value = 0.1977100304534503
$ws.Range("K3").Value = 0.9563608868237736
$ws.Range("B4").Value = 0.1273558660456899
$ws.Range("D4").Value = 0.1209730199446852
$ws.Range("E4").Value = 0.1438026769056862
$ws.Range("F4").Value = 2.029533379137362
$ws.Range("G4").Value = 1.365185102181968
$ws.Range("H4").Value = 1.269937368049696
$ws.Range("J4").Value = 0.1952319234243163
$ws.Range("K4").Value = 0.893245966995778
$ws.Range("B5").Value = 0.1250724214932291
$ws.Range("D5").Value = 0.1201559867123478
$ws.Range("E5").Value = 0.1429940104873602
$ws.Range("F5").Value = 2.028691802289856
$ws.Range("G5").Value = 1.364343669667562
$ws.Range("H5").Value = 1.271303719478766
$ws.Range("J5").Value = 0.1942525716327737
$ws.Range("K5").Value = 0.8675807508392381
$ws.Range("B6").Value = 0.1246943666990745
$ws.Range("D6").Value = 0.1200209689938063
$ws.Range("E6").Value = 0.1428608744245921
$ws.Range("F6").Value = 2.028571379293382
$ws.Range("G6").Value = 1.364218979579718
$ws.Range("H6").Value = 1.271541080968433
$ws.Range("J6").Value = 0.1940917909191526
$ws.Range("K6").Value = 0.8633223758525617
$ws.Range("B7").Value = 0.1273249964113035
$ws.Range("D7").Value = 0.1209619576123089
$ws.Range("E7").Value = 0.1437916943640509
$ws.Range("F7").Value = 2.029520733704501
$ws.Range("G7").Value = 1.365172746262701
$ws.Range("H7").Value = 1.269955092435467
$ws.Range("J7").Value = 0.1952185921544043
$ws.Range("K7").Value = 0.8928996150539774
$ws.Range("B8").Value = 0.1390938771850898
$ws.Range("D8").Value = 0.1252083888545172
$ws.Range("E8").Value = 0.1480652503470132
$ws.Range("F8").Value = 2.036523242811342
$ws.Range("G8").Value = 1.371572016876328
$ws.Range("H8").Value = 1.264481564447209
$ws.Range("J8").Value = 0.2004588455962732
$ws.Range("K8").Value = 1.023895710135889
$ws.Range("B9").Value = 0.1627434781141233
$ws.Range("D9").Value = 0.1338702013345738
$ws.Range("E9").Value = 0.1570431944190105
$ws.Range("F9").Value = 2.060417525360606
$ws.Range("G9").Value = 1.392054049768177
$ws.Range("H9").Value = 1.259247833233559
$ws.Range("J9").Value = 0.2117034941033751
$ws.Range("K9").Value = 1.282479110133863
$ws.Range("B10").Value = 0.1804640725651723
$ws.Range("D10").Value = 0.1404359668794655
$ws.Range("E10").Value = 0.1640054435946467
$ws.Range("F10").Value = 2.084222417539209
$ws.Range("G10").Value = 1.411999279165087
$ws.Range("H10").Value = 1.258774748130548
$ws.Range("J10").Value = 0.2205621039800576
$ws.Range("K10").Value = 1.473508568358454
$ws.Range("B11").Value = 0.1886001336739582
$ws.Range("D11").Value = 0.1434663899115947
$ws.Range("E11").Value = 0.1672527382669671
$ws.Range("F11").Value = 2.09642109097588
$ws.Range("G11").Value = 1.422149612809136
$ws.Range("H11").Value = 1.259296061740372
$ws.Range("J11").Value = 0.2247233203064667
$ws.Range("K11").Value = 1.560644644420677
$ws.Range("B12").Value = 0.1916917329911314
$ws.Range("D12").Value = 0.1446201596995422
$ws.Range("E12").Value = 0.1684939483322907
$ws.Range("F12").Value = 2.101238266333837
$ws.Range("G12").Value = 1.426149187778066
$ws.Range("H12").Value = 1.259599708212107
$ws.Range("J12").Value = 0.2263180602374888
$ws.Range("K12").Value = 1.593674645231033
$ws.Range("B13").Value = 0.1910254299825596
$ws.Range("D13").Value = 0.144371399269275
$ws.Range("E13").Value = 0.1682261186123668
$ws.Range("F13").Value = 2.100191990677516
$ws.Range("G13").Value = 1.425280860152526
$ws.Range("H13").Value = 1.259529582411858
$ws.Range("J13").Value = 0.2259737590780304
$ws.Range("K13").Value = 1.586559562333377
$ws.Range("B14").Value = 0.1888542686470771
$ws.Range("D14").Value = 0.1435611870327591
$ws.Range("E14").Value = 0.1673546222428541
$ws.Range("F14").Value = 2.096813433806176
$ws.Range("G14").Value = 1.422475530928949
$ws.Range("H14").Value = 1.259318911912544
$ws.Range("J14").Value = 0.2248541396966317
$ws.Range("K14").Value = 1.563361374368128
$ws.Range("B15").Value = 0.1875257530533361
$ws.Range("D15").Value = 0.1430657162167535
$ws.Range("E15").Value = 0.1668223073310315
$ws.Range("F15").Value = 2.094769754907091
$ws.Range("G15").Value = 1.420777512742944
$ws.Range("H15").Value = 1.259203714950871
$ws.Range("J15").Value = 0.224170814450801
$ws.Range("K15").Value = 1.549156158944868
$ws.Range("B16").Value = 0.1799338584611405
$ws.Range("D16").Value = 0.1402387943827392
$ws.Range("E16").Value = 0.1637948383080428
$ws.Range("F16").Value = 2.083452839273235
$ws.Range("G16").Value = 1.411357688920987
$ws.Range("H16").Value = 1.258755527521714
$ws.Range("J16").Value = 0.2202928097580354
$ws.Range("K16").Value = 1.467818727030703
$ws.Range("B17").Value = 0.1752955704709365
$ws.Range("D17").Value = 0.1385157003121549
$ws.Range("E17").Value = 0.1619581127508667
$ws.Range("F17").Value = 2.07686168940586
$ws.Range("G17").Value = 1.405855505696309
$ws.Range("H17").Value = 1.258669450733464
$ws.Range("J17").Value = 0.2179474910036419
$ws.Range("K17").Value = 1.4179809170368
$ws.Range("B18").Value = 0.1726348043534074
$ws.Range("D18").Value = 0.1375287316472793
$ws.Range("E18").Value = 0.1609092187711099
$ws.Range("F18").Value = 2.073199507298995
$ws.Range("G18").Value = 1.402792123534226
$ws.Range("H18").Value = 1.258689247033885
$ws.Range("J18").Value = 0.2166108876818811
$ws.Range("K18").Value = 1.389337748439402
$ws.Range("B19").Value = 0.1717351298534737
$ws.Range("D19").Value = 0.1371952686262574
$ws.Range("E19").Value = 0.1605553764336634
$ws.Range("F19").Value = 2.071981663890213
$ws.Range("G19").Value = 1.401772288620919
$ws.Range("H19").Value = 1.25870784310942
$ws.Range("J19").Value = 0.216160457410453
$ws.Range("K19").Value = 1.379643494426773
$ws.Range("B20").Value = 0.1757885948855744
$ws.Range("D20").Value = 0.1386987017307746
$ws.Range("E20").Value = 0.1621528547065978
$ws.Range("F20").Value = 2.077549984663278
$ws.Range("G20").Value = 1.406430728420105
$ws.Range("H20").Value = 1.25867143858602
$ws.Range("J20").Value = 0.2181958741644792
$ws.Range("K20").Value = 1.42328394048036
$ws.Range("B21").Value = 0.1894917032676204
$ws.Range("D21").Value = 0.143798997812425
$ws.Range("E21").Value = 0.1676102888689144
$ws.Range("F21").Value = 2.097800422833444
$ws.Range("G21").Value = 1.423295286263993
$ws.Range("H21").Value = 1.25937790511577
$ws.Range("J21").Value = 0.2251824833282114
$ws.Range("K21").Value = 1.570174339876019
$ws.Range("B22").Value = 0.1985094963958005
$ws.Range("D22").Value = 0.147168531753934
$ws.Range("E22").Value = 0.171244247465026
$ws.Range("F22").Value = 2.112188592585099
$ws.Range("G22").Value = 1.435226197173108
$ws.Range("H22").Value = 1.260459001259079
$ws.Range("J22").Value = 0.229859292101068
$ws.Range("K22").Value = 1.666370546458495
$ws.Range("B23").Value = 0.1936908942986975
$ws.Range("D23").Value = 0.1453668552494776
$ws.Range("E23").Value = 0.1692985830822238
$ws.Range("F23").Value = 2.104403546224859
$ws.Range("G23").Value = 1.428774954901854
$ws.Range("H23").Value = 1.259825218493802
$ws.Range("J23").Value = 0.2273530382251892
$ws.Range("K23").Value = 1.61501113494711
$ws.Range("B24").Value = 0.1755656803253629
$ws.Range("D24").Value = 0.1386159553819084
$ws.Range("E24").Value = 0.1620647898664629
$ws.Range("F24").Value = 2.077238410356316
$ws.Range("G24").Value = 1.406170359133569
$ws.Range("H24").Value = 1.258670324089252
$ws.Range("J24").Value = 0.2180835435402884
$ws.Range("K24").Value = 1.420886414731285
$ws.Range("B25").Value = 0.1562848090735827
$ws.Range("D25").Value = 0.1314913476145989
$ws.Range("E25").Value = 0.1545502623181818
$ws.Range("F25").Value = 2.052859918908155
$ws.Range("G25").Value = 1.385657546078264
$ws.Range("H25").Value = 1.260072915470275
$ws.Range("J25").Value = 0.2085571397841335
$ws.Range("K25").Value = 1.21234230113015
